$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.270.55"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "3.139.77"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.71%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  -5.95%  "

$ws.Range("D9").Value = "3.153.72"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("E10").Value = "  -3.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.59"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.384"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "3.695.84"
$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "64.295.19"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.08"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("D17").Value = "3.149.68"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("E18").Value = "  -3.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "400.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.59%  "

$ws.Range("E21").Value = "  -3.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.40"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.482"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.196"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.45%  "

$ws.Range("E27").Value = "  -4.99%  "

$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.79"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.82"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.96%  "

$ws.Range("E36").Value = "  -2.95%  "

$ws.Range("E37").Value = "  -2.58%  "

$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("D39").Value = "2.645.70"
$ws.Range("E39").Value = "  -3.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.05"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.689"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0610"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.40%  "

$ws.Range("E46").Value = "  -4.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "285.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.996"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0971"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("E51").Value = "  +0.00%  "
